$wb = $excel.ActiveWorkbook

$uk = $wb.Worksheets.Item("UK")
$belgium = $wb.Worksheets.Item("Belgium")

# --- Create "Romania" as a copy of the "UK" sheet, placed after "Belgium" ---
$uk.Copy($null, $belgium)
$romania = $wb.Worksheets.Item($belgium.Index + 1)
$romania.Name = "Romania"

# --- Create "Slovakia" as a copy of the "UK" sheet, placed after "Romania" ---
$uk.Copy($null, $romania)
$slovakia = $wb.Worksheets.Item($romania.Index + 1)
$slovakia.Name = "Slovakia"

# --- Populate Romania's market info ---
$romania.Rows.Item(2).AutoFit() | Out-Null
$romania.Range("B3").Copy() | Out-Null
$romania.Range("B4").PasteSpecial(-4122) | Out-Null
$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3534/T3544"

# --- Populate Slovakia's market info ---
$slovakia.Rows.Item(2).AutoFit() | Out-Null
$slovakia.Range("B2").Value = "Slovakia market"
$slovakia.Range("B4").Value = "NGC-4306/T3560/T3567"

# --- Selections / active sheet bookkeeping ---
$uk.Activate()
$uk.Range("A1:XFD1048576").Select() | Out-Null

$romania.Activate()
$romania.Range("B2:B4").Select() | Out-Null

$slovakia.Activate()
$slovakia.Range("B5").Select() | Out-Null
